# Adds the Artisan Command `replayLookahead`
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# Insert a new row after the current "pidLookahead" row (row 113) so everything
# below shifts down by one, then fill in the new row contents.
$ws.Rows.Item(114).Insert()

# update text for the existing pidLookahead description (row 113)
$ws.Range("C113").Value = "sets the PID lookahead in seconds"

# populate the newly inserted row 114 with the new command
$ws.Range("B114").Value = "replayLookahead(<int>)"
$ws.Range("C114").Value = "sets the Ramping Event Replay lookahead in seconds"

# re-assign these two cells as plain text (collapsing their previous multi-run
# rich text formatting down to a single uniform, non-italic string)
$ws.Range("B28").Value = "mwrite(deviceID,register,andMask,orMask) or mwrite(deviceID,register,andMask,orMask,value)"
$ws.Range("B30").Value = "writeBCD(deviceID,register,value) or writeBCD([deviceID,register,value],..,[deviceID,register,value])"
